$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion message text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 11.4 = 46470.52 pesos`n✅ 46470.52 pesos = 11.34 = 965.24 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $text

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 87.69
$ws2.Range("O10").Value = 4075
$ws2.Range("N12").Value = 4099
$ws2.Range("O12").Value = 85.14
